$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row of data (row 22) mirroring the style/pattern of row 21
$ws.Range("A22").Value = "Danh project"
$ws.Range("B22").Value = "Make this damn thing"
$ws.Range("C22").Value = "Hung"
$ws.Range("D22").Value = "In progress"
$ws.Range("E22").Value = "1st draft sent and approved on 28 Oct"
$ws.Range("F22").Value = "S"

# Match formatting used by the row above it (A21 bold, B21:F21 regular)
$ws.Range("A21").Copy()
$ws.Range("A22").PasteSpecial(-4122)
$ws.Range("B21:F21").Copy()
$ws.Range("B22:F22").PasteSpecial(-4122)

# Update selection to reflect where the cursor ended up after editing
$ws.Range("E27").Select()
